# Reparando errores de conexion
# Adds a new "hallazgo" (finding) row (row 13) to the Hoja1 worksheet,
# mirroring the formatting used by the existing rows (B column uses the
# "alert" style from rows 11/12, C:G use the plain wrap-text style used
# throughout column G), then fills in the new finding's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed row 13 with formatting copied from existing cells so the new row
# re-uses the workbook's existing style entries instead of minting new
# ones: B13 copies B12's style (bold red "Alta Criticidad"-style font with
# border), C13:G13 copy G2's style (plain wrap-text, no border) which is
# the same style already used for every cell in column G.
$ws.Range("B12").Copy($ws.Range("B13"))
$ws.Range("G2").Copy($ws.Range("C13:G13"))

# Fill in the new finding's content.
$ws.Range("B13").Value = "Retirar consultar inscrustadas en el html, nivel critico de seguridad"
$ws.Range("C13").Value = "Medio"
$ws.Range("D13").Value = "incrustacion de sentencias sql"
$ws.Range("E13").Value = "Un usuario con un nivel basico o intermedio de conocimiento de programacion puede ver y acceder  a la al codigo fuente y ver incrustacion de sentencias/ mala practica/ vulnerabilidad en la Pagina"
$ws.Range("F13").Value = "Error de codigo fuente - innefecciencia en la seguridad de la informacion"
$ws.Range("G13").Value = "Abierta"

# Match the row height Excel computed for the wrapped text in the new row.
$ws.Rows.Item(13).RowHeight = 78.75

# Move the selection the way the author left it after typing the new row.
$ws.Range("G15").Select() | Out-Null
